$d = $word.ActiveDocument

# Helper: find $searchText inside $range (without performing a Find-driven
# replace, which would trigger smart-quote autocorrect on the replacement
# text, turning straight quotes/apostrophes into curly ones) and then
# overwrite just the matched span via Range.Text, which preserves the
# original text verbatim (including straight apostrophes).
function Replace-InRange($range, $searchText, $replaceText) {
    $rng = $range.Duplicate
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $doc = $rng.Document
        $target = $doc.Range($rng.Start, $rng.End)
        $target.Text = $replaceText
        return $true
    }
    return $false
}

# --- "24 лютого 2021, 17:00, конференція в Zoom" -> "... березня ..." ------
$p9 = $d.Paragraphs.Item(9)
if (-not (Replace-InRange $p9.Range " лютого 20" " березня 20")) {
    Replace-InRange $d.Content " лютого 20" " березня 20" | Out-Null
}

# --- Speaker name: "Богачов Кирило" -> "Ткачев Ігор Іванович" --------------
$p12 = $d.Paragraphs.Item(12)
if (-not (Replace-InRange $p12.Range "Богачов Кирило" "Ткачев Ігор Іванович")) {
    Replace-InRange $d.Content "Богачов Кирило" "Ткачев Ігор Іванович" | Out-Null
}

# --- Speaker affiliation -----------------------------------------------------
$p13 = $d.Paragraphs.Item(13)
$oldAffil = "(студент Фізико-технічного інституту НТУ України «Київський політехнічний університет імені Ігоря Сікорського»)"
$newAffil = "(к.ф.-м.н., Науково-навчальний центр прикладної інформатики НАН України)"
if (-not (Replace-InRange $p13.Range $oldAffil $newAffil)) {
    Replace-InRange $d.Content $oldAffil $newAffil | Out-Null
}

# --- Talk title ---------------------------------------------------------------
$p15 = $d.Paragraphs.Item(15)
$oldTitle = "Аналіз фотоемісійних спектрів багатозонних надпровідників методами машинного навчання з застосуванням згорткових нейронних мереж"
$newTitle = "Класична задача дискретної математики — мінімізація диз'юнктивної нормальної форми у формалізмі реляційних схем"
if (-not (Replace-InRange $p15.Range $oldTitle $newTitle)) {
    Replace-InRange $d.Content $oldTitle $newTitle | Out-Null
}

# --- "Середа, 24 лютого 2021р. ... Zoom" -> "... березня ..." --------------
$p28 = $d.Paragraphs.Item(28)
if (-not (Replace-InRange $p28.Range "лютого" "березня")) {
    Replace-InRange $d.Content "лютого" "березня" | Out-Null
}

# Remove one redundant tab (13 -> 12) right after "р." in that same
# paragraph (re-fetched, since the text above may have shifted offsets).
$p28 = $d.Paragraphs.Item(28)
$t = $p28.Range.Text
$tabIdx = $t.IndexOf([char]9)
if ($tabIdx -ge 0) {
    $absStart = $p28.Range.Start + $tabIdx
    $tabRange = $d.Range($absStart, $absStart + 1)
    $tabRange.Delete()
}
